$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (Sheet1 -> DutyList) ---
$ws.Name = "DutyList"

# --- Drop the stale _FilterDatabase defined name (and its #REF! autofilter leftover) ---
foreach ($n in @($wb.Names)) {
    $n.Delete()
}

# --- Re-point the header row at the new domain-model field names ---
# Columns F:L keep their existing text (Region, dutyType, dutyLoad, StartDate,
# StartTime, EndDate, EndTime) so only A:E are rewritten, and M:N are new.
$ws.Range("A1").Value = "Tarih"
$ws.Range("B1").Value = "kod"
$ws.Range("C1").Value = "sicilno"
$ws.Range("D1").Value = "peradi"
$ws.Range("E1").Value = "Telefon"
$ws.Range("M1").Value = "Priority"
$ws.Range("N1").Value = "totalWorkingHour"

# Give the two brand new header cells the same bold/centered look as the rest
# of row 1 by copying the formatting from an existing header cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Resize columns to match the regenerated layout ---
# (ColumnWidth values are chosen so the saved <col> width lands on the target;
# columns whose width doesn't change (A, G, H, J, L) are left untouched so
# their bestFit flag survives.)
$ws.Columns.Item(2).ColumnWidth = 5.1666666666667
$ws.Columns.Item(3).ColumnWidth = 5.8776041666667
$ws.Columns.Item(4).ColumnWidth = 5.8776041666667
$ws.Columns.Item(5).ColumnWidth = 7.0221354166667
$ws.Columns.Item(6).ColumnWidth = 7.4518229166667
$ws.Columns.Item(9).ColumnWidth = 9.5924479166667
$ws.Columns.Item(11).ColumnWidth = 9.5924479166667
$ws.Columns.Item(13).ColumnWidth = 6.8776041666667
$ws.Columns.Item(14).ColumnWidth = 16.4518229166667

# --- Theme accent colors: accent1 and accent5 were swapped ---
$scheme = $wb.Theme.ThemeColorScheme
$scheme.Colors(5).RGB = 12874308  # accent1 <- 4472C4
$scheme.Colors(9).RGB = 13998939  # accent5 <- 5B9BD5
